$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Running suites c and d: flip Runmode from N to Y for the remaining
# AuthoringTest-family rows so they are included in the run.
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"

# Active cell/selection moved to C9 on this sheet.
$ws.Range("C9").Select()
